$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 210-211 (formatting copied from row above, matching
# the existing date format used by column D in the surrounding rows).
$ws.Rows("210:211").Insert()

# New row 210
$ws.Range("A210").Value2 = 10
$ws.Range("B210").Value2 = "Vega Modelo de Temuco"
$ws.Range("C210").Value2 = "La Araucanía"
$ws.Range("D210").Value2 = 44767
$ws.Range("E210").Value2 = 9
$ws.Range("F210").Value2 = 100112017
$ws.Range("G210").Value2 = "Apio"
$ws.Range("H210").Value2 = "Americana (o)"
$ws.Range("I210").Value2 = "Primera"
$ws.Range("J210").Value2 = 200
$ws.Range("K210").Value2 = 13000
$ws.Range("L210").Value2 = 13000
$ws.Range("M210").Value2 = 13000
$ws.Range("N210").Value2 = "$/docena de matas"
$ws.Range("O210").Value2 = "Provincia del Elquí"
$ws.Range("P210").Value2 = 2167
$ws.Range("Q210").Value2 = 6
$ws.Range("R210").Value2 = "Hortaliza"

# New row 211
$ws.Range("A211").Value2 = 10
$ws.Range("B211").Value2 = "Vega Modelo de Temuco"
$ws.Range("C211").Value2 = "La Araucanía"
$ws.Range("D211").Value2 = 44767
$ws.Range("E211").Value2 = 9
$ws.Range("F211").Value2 = 100112017
$ws.Range("G211").Value2 = "Apio"
$ws.Range("H211").Value2 = "Americana (o)"
$ws.Range("I211").Value2 = "Segunda"
$ws.Range("J211").Value2 = 80
$ws.Range("K211").Value2 = 10000
$ws.Range("L211").Value2 = 10000
$ws.Range("M211").Value2 = 10000
$ws.Range("N211").Value2 = "$/docena de matas"
$ws.Range("O211").Value2 = "Provincia del Elquí"
$ws.Range("P211").Value2 = 1667
$ws.Range("Q211").Value2 = 6
$ws.Range("R211").Value2 = "Hortaliza"
